# Replace the hard-coded "www.drpaulduenas.com" footer text with a
# configurable MERGEFIELD ("website") field, matching the pattern already
# used by the other footer fields (branch office address/phone, emergency
# number, city, etc.).

$d = $word.ActiveDocument

# The text lives in the primary footer (Section 1, Footer index 1).
$footer = $d.Sections.First.Footers.Item(1)

# Locate the exact run that needs to change.
$rng = $footer.Range
$found = $rng.Find.Execute("www.drpaulduenas.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Clear the literal text but keep the run / paragraph (and its
    # formatting, jc=center, bold Avenir Book 20, etc.) intact.
    $rng.Text = ""

    # Re-insert the paragraph with the MERGEFIELD field construct
    # (begin / instrText / separate / result / end), reusing the same
    # paragraph identity and run formatting as the original.
    $xml = '<w:p w14:paraId="24EA949D" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRDefault="004D2A29" w:rsidP="004D2A29" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«=website»</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>'

    $rng.InsertXML($xml)
}

Write-Output "done"
